$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Roll the "twelve months ended" year headers forward by one year:
#    drop 1396/12, keep 1397..1400, add 1401/12 as the new trailing column.
#    (Columns E..I = the 5 year columns in the two header rows 8 and 24.)
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 5).Value  = "دوازده ماهه منتهی به 1397/12"
$ws.Cells.Item(8, 6).Value  = "دوازده ماهه منتهی به 1398/12"
$ws.Cells.Item(8, 7).Value  = "دوازده ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 8).Value  = "دوازده ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 9).Value  = "دوازده ماهه منتهی به 1401/12"

$ws.Cells.Item(24, 5).Value = "دوازده ماهه منتهی به 1397/12"
$ws.Cells.Item(24, 6).Value = "دوازده ماهه منتهی به 1398/12"
$ws.Cells.Item(24, 7).Value = "دوازده ماهه منتهی به 1399/12"
$ws.Cells.Item(24, 8).Value = "دوازده ماهه منتهی به 1400/12"
$ws.Cells.Item(24, 9).Value = "دوازده ماهه منتهی به 1401/12"

# ---------------------------------------------------------------------------
# 2) Update the database: every yearly figure shifts one column to the left
#    (the oldest year drops off) and a freshly computed value lands in the
#    new trailing (1401) column. A handful of cells also got a small
#    "read_price algorithm" correction versus a pure shift.
# ---------------------------------------------------------------------------

# هزینه حمل و نقل و انتقال
$ws.Cells.Item(10, 5).Value = 450
$ws.Cells.Item(10, 6).Value = 317
$ws.Cells.Item(10, 7).Value = 456
$ws.Cells.Item(10, 8).Value = 237
$ws.Cells.Item(10, 9).Value = 0

# هزینه خدمات پس از فروش
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0

# حق العمل و کمیسیون فروش
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0

# هزینه تبلیغات
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0

# هزینه مواد مصرفی
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0

# هزینه انرژی (آب، برق، گاز و سوخت)
$ws.Cells.Item(15, 5).Value = 271
$ws.Cells.Item(15, 6).Value = 311
$ws.Cells.Item(15, 7).Value = 311
$ws.Cells.Item(15, 8).Value = 342
$ws.Cells.Item(15, 9).Value = 492

# هزینه استهلاک
$ws.Cells.Item(16, 5).Value = 3080
$ws.Cells.Item(16, 6).Value = 3096
$ws.Cells.Item(16, 7).Value = 3088
$ws.Cells.Item(16, 8).Value = 5848
$ws.Cells.Item(16, 9).Value = 8274

# هزینه حقوق و دستمزد
$ws.Cells.Item(17, 5).Value = 16577
$ws.Cells.Item(17, 6).Value = 23037
$ws.Cells.Item(17, 7).Value = 31032
$ws.Cells.Item(17, 8).Value = 121930
$ws.Cells.Item(17, 9).Value = 186135

# هزینه مطالبات مشکوک الوصول
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0

# سایر هزینه ها
$ws.Cells.Item(19, 5).Value = 46241
$ws.Cells.Item(19, 6).Value = 57450
$ws.Cells.Item(19, 7).Value = 64954
$ws.Cells.Item(19, 8).Value = 49593
$ws.Cells.Item(19, 9).Value = 70729

# جمع
$ws.Cells.Item(20, 5).Value = 66619
$ws.Cells.Item(20, 6).Value = 84211
$ws.Cells.Item(20, 7).Value = 99841
$ws.Cells.Item(20, 8).Value = 177950
$ws.Cells.Item(20, 9).Value = 265630

# تعداد پرسنل غیر تولیدی شرکت
$ws.Cells.Item(26, 5).Value = 100
$ws.Cells.Item(26, 6).Value = 97
$ws.Cells.Item(26, 7).Value = 95
$ws.Cells.Item(26, 8).Value = 93
$ws.Cells.Item(26, 9).Value = 93

# تعداد پرسنل تولیدی شرکت
$ws.Cells.Item(27, 5).Value = 315
$ws.Cells.Item(27, 6).Value = 300
$ws.Cells.Item(27, 7).Value = 325
$ws.Cells.Item(27, 8).Value = 342
$ws.Cells.Item(27, 9).Value = 334
